# Add a new "ES Bribery Incidence" column (U) to the Datasets and Years
# tracker sheet, with checkmarks for years 2006-2016 and the "Years Vary
# by Country" summary row (rows 20-31), matching the style already used
# for the other header cells / checkmark cells in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell U3: "ES Bribery Incidence" -------------------------------
# Clone the formatting of an existing header cell (B3: Times New Roman 12,
# automatic/theme color) so we reuse the same cellXfs/font instead of
# inventing a brand-new, differently-ordered style entry, then tweak the
# font color explicitly (matches the author's source file, which carries
# an explicit theme color rather than "no color element").
$ws.Range("B3").Copy()
$ws.Range("U3").PasteSpecial(-4122)
$ws.Range("U3").Value = "ES Bribery Incidence"
$ws.Range("U3").Font.ThemeColor = 1

# --- Checkmarks for U20:U31 -------------------------------------------------
# Clone the formatting of an existing checkmark cell (B4) first so the new
# cells pick up the same style index as the rest of the checkmark column,
# then fill in the checkmark glyph used throughout the sheet.
$ws.Range("B4").Copy()
$ws.Range("U20:U31").PasteSpecial(-4122)
$ws.Range("U20:U31").Value = "✓"

# Page was set to print in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it.
$ws.Range("U31").Select() | Out-Null
